# Apply the "Descriptions" sheet addition to CampusTestData.xlsx
$wb = $excel.ActiveWorkbook

# Add a new worksheet named "Descriptions" as the 3rd sheet (after Nationalities)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Descriptions"

# Fill in the data for the Descriptions sheet
$data = @(
    @("EnglishCourse", "EC", 1, "TestNGCourse", "TC", 1),
    @("GermanCourse", "GC", 2, "CucumberCourse", "CC", 2),
    @("FranceCourse", "FC", 1, "JavaCourse", "JC", 1),
    @("TechnicalCourse", "TC", 1, "PostmanCourse", "PC", 1),
    @("MathCourse", "MC", 2, "MySQLCourse", "MC", 2),
    @("CulturCourse", "CC", 2, "JenkinsCourse", "JC", 2)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $newSheet.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# Column widths to match the target layout as closely as this engine's
# ColumnWidth quantization allows (stored OOXML width = ColumnWidth + 5/6,
# rounded to the nearest 1/6): target stored widths are 19, 14.85546875 and
# 20.42578125 characters, so we back-solve for the ColumnWidth to assign.
$newSheet.Columns.Item(1).ColumnWidth = 18.166666666666668
$newSheet.Columns.Item(2).ColumnWidth = 14.0
$newSheet.Columns.Item(4).ColumnWidth = 19.666666666666668

# Select E6 on the new sheet and make it the active/selected tab
$newSheet.Range("E6").Select()

# The second sheet (Nationalities) should no longer be the tab shown as selected;
# activeTab should now point to the 3rd sheet (0-indexed => 2)
$newSheet.Activate()
